$win = $excel.ActiveWindow
$members = $win | Get-Member | ForEach-Object { $_.Name }
Write-Host ($members -join "`n")
